$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the additional row 2 values
$ws.Range("B2").Value = 44
$ws.Range("C2").Value = 324
$ws.Range("D2").Value = "rrg"
$ws.Range("E2").Value = "g4"
$ws.Range("G2").Value = 25
$ws.Range("I2").Value = "wer"

# Add row 10 with concatenation formulas for each column A..I
$ws.Range("A10").Formula = "=A1&A2"
$ws.Range("B10").Formula = "=B1&B2"
$ws.Range("C10").Formula = "=C1&C2"
$ws.Range("D10").Formula = "=D1&D2"
$ws.Range("E10").Formula = "=E1&E2"
$ws.Range("F10").Formula = "=F1&F2"
$ws.Range("G10").Formula = "=G1&G2"
$ws.Range("H10").Formula = "=H1&H2"
$ws.Range("I10").Formula = "=I1&I2"

# Update the active selection to match the target (row 10 selected)
$ws.Range("A10:I10").Select()
